# Update "想去人数" (want-to-go count) values for a few conventions.
# Sheet "展览" (exhibitions) - rows 2, 4, 5
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 253
$ws1.Range("F4").Value = 838
$ws1.Range("F5").Value = 525

# Sheet "全部类型" (all types) - rows 2, 4, 6
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 253
$ws4.Range("F4").Value = 838
$ws4.Range("F6").Value = 525
